# Social network scale added
#
# 1. "socialnetwork" sheet (sheet3.xml): remove the "meet/speak/write"
#    frequency rows for each wave (child/family/friend contact-mode rows),
#    keeping only the social-network-scale rows (spouse/child/family/friends
#    + "closeness" rows). These live in five contiguous 9-row blocks.
# 2. Update the selections on a couple of sheets and the active sheet.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("socialnetwork")

# Delete the obsolete rows bottom-to-top so earlier row numbers stay valid.
$ws3.Rows("78:86").Delete()
$ws3.Rows("61:69").Delete()
$ws3.Rows("44:52").Delete()
$ws3.Rows("27:35").Delete()
$ws3.Rows("10:18").Delete()

# lifesatisfaction: move selection to D2 (was D14)
$ws2 = $wb.Worksheets.Item("lifesatisfaction")
[void]$ws2.Range("D2").Select()

# socialnetwork: move selection to A42 (was the whole-column A1:D1048576)
[void]$ws3.Range("A42").Select()

# socialnetwork becomes the active tab (was demographics)
[void]$ws3.Activate()
